{"js": "const body = context.document.body;\n\nasync function replaceOnce(oldText, newText) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n\nawait replaceOnce(\n  \"n the year a particular Programming Language has appeared on GitHub\",\n  \"n the year a particular Programming Language was released\"\n);\n\nawait replaceOnce(\n  \"nature of programming languages?\",\n  \"nature of the Programming Languages used in GitHub repositories?\"\n);\n\nawait replaceOnce(\n  \" What percentage of those clusters have appeared on GitHub prior to 2000?\",\n  \" What percentage of those clusters contain Programming Languages that were released prior to 1993?\"\n);\n\nawait replaceOnce(\n  \" determines to have a coefficient, create another Linear Regression model. How did this effect the original model? What does this tell you about the data?\",\n  \" determines to have a coefficient, create another Linear Regression model. How did this effect the original model? What does this tell you about the GitHut data regarding Programming Languages used in GitHub repositories?\"\n);\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{\n        Find = \"n the year a particular Programming Language has appeared on GitHub\"\n        Replace = \"n the year a particular Programming Language was released\"\n    },\n    @{\n        Find = \"nature of programming languages?\"\n        Replace = \"nature of the Programming Languages used in GitHub repositories?\"\n    },\n    @{\n        Find = \" What percentage of those clusters have appeared on GitHub prior to 2000?\"\n        Replace = \" What percentage of those clusters contain Programming Languages that were released prior to 1993?\"\n    },\n    @{\n        Find = \" determines to have a coefficient, create another Linear Regression model. How did this effect the original model? What does this tell you about the data?\"\n        Replace = \" determines to have a coefficient, create another Linear Regression model. How did this effect the original model? What does this tell you about the GitHut data regarding Programming Languages used in GitHub repositories?\"\n    }\n)\n\nforeach ($r in $replacements) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $r.Find\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $r.Replace\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.Execute(\n        $r.Find,\n        $true,\n        $false,\n        $false,\n        $false,\n        $false,\n        $true,\n        1,\n        $false,\n        $r.Replace,\n        2\n    )\n}\n"}
